$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column (D) retains text formatting for numeric-looking strings
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "67.869.79"
$ws.Range("E2").Value = "  +1.75%  "

# Row 3
$ws.Range("D3").Value = "3.340.00"
$ws.Range("E3").Value = "  +2.58%  "

# Row 4
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.16%  "

# Row 5
$ws.Range("D5").Value = "583.46"
$ws.Range("E5").Value = "  +2.91%  "

# Row 6
$ws.Range("D6").Value = "176.90"
$ws.Range("E6").Value = "  +1.62%  "

# Row 7
$ws.Range("E7").Value = "  -0.28%  "

# Row 8
$ws.Range("E8").Value = "  +2.05%  "

# Row 9
$ws.Range("D9").Value = "3.335.44"
$ws.Range("E9").Value = "  +2.76%  "

# Row 10
$ws.Range("D10").Value = "0.182"
$ws.Range("E10").Value = "  +5.71%  "

# Row 11
$ws.Range("E11").Value = "  +2.77%  "

# Row 12
$ws.Range("D12").Value = "46.72"
$ws.Range("E12").Value = "  +3.90%  "

# Row 13
$ws.Range("D13").Value = "0.0000273"
$ws.Range("E13").Value = "  +2.64%  "

# Row 14
$ws.Range("D14").Value = "691.62"
$ws.Range("E14").Value = "  +0.18%  "

# Row 15
$ws.Range("D15").Value = "3.884.06"
$ws.Range("E15").Value = "  +2.47%  "

# Row 16
$ws.Range("E16").Value = "  +2.70%  "

# Row 17
$ws.Range("D17").Value = "67.907.95"
$ws.Range("E17").Value = "  +1.13%  "

# Row 18
$ws.Range("E18").Value = "  -0.12%  "

# Row 19
$ws.Range("D19").Value = "3.336.46"
$ws.Range("E19").Value = "  +1.69%  "

# Row 20
$ws.Range("D20").Value = "17.42"
$ws.Range("E20").Value = "  +1.49%  "

# Row 21
$ws.Range("D21").Value = "11.11"
$ws.Range("E21").Value = "  +4.75%  "

# Row 22
$ws.Range("D22").Value = "0.896"
$ws.Range("E22").Value = "  +1.89%  "

# Row 23
$ws.Range("D23").Value = "5.38"
$ws.Range("E23").Value = "  +6.03%  "

# Row 24
$ws.Range("D24").Value = "17.08"
$ws.Range("E24").Value = "  +1.50%  "

# Row 25
$ws.Range("D25").Value = "98.63"
$ws.Range("E25").Value = "  +1.45%  "

# Row 26
$ws.Range("D26").Value = "3.89"
$ws.Range("E26").Value = "  +1.18%  "

# Row 27
$ws.Range("D27").Value = "2.70"
$ws.Range("E27").Value = "  +0.82%  "

# Row 28
$ws.Range("D28").Value = "9.52"
$ws.Range("E28").Value = "  +3.62%  "

# Row 29
$ws.Range("D29").Value = "32.97"
$ws.Range("E29").Value = "  +1.63%  "

# Row 30
$ws.Range("D30").Value = "8.58"
$ws.Range("E30").Value = "  +3.17%  "

# Row 31
$ws.Range("D31").Value = "7.13"
$ws.Range("E31").Value = "  +7.73%  "

# Row 32
$ws.Range("D32").Value = "572.93"
$ws.Range("E32").Value = "  -0.47%  "

# Row 33
$ws.Range("D33").Value = "11.01"
$ws.Range("E33").Value = "  +3.11%  "

# Row 34
$ws.Range("D34").Value = "0.106"
$ws.Range("E34").Value = "  +3.45%  "

# Row 35
$ws.Range("B35").Value = "Maker"
$ws.Range("C35").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D35").Value = "3.723.46"
$ws.Range("E35").Value = "  -2.82%  "

# Row 36
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").Value = "57.31"
$ws.Range("E36").Value = "  +3.96%  "

# Row 37
$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.11%  "

# Row 38
$ws.Range("D38").Value = "3.37"
$ws.Range("E38").Value = "  +3.95%  "

# Row 39
$ws.Range("B39").Value = "InjectiveProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D39").Value = "34.15"
$ws.Range("E39").Value = "  +8.63%  "

# Row 40
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "0.131"
$ws.Range("E40").Value = "  +2.73%  "

# Row 41
$ws.Range("D41").Value = "3.20"
$ws.Range("E41").Value = "  +4.86%  "

# Row 42
$ws.Range("D42").Value = "2.65"
$ws.Range("E42").Value = "  +3.23%  "

# Row 43
$ws.Range("D43").Value = "0.0₃0675"
$ws.Range("E43").Value = "  +1.90%  "

# Row 44
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").Value = "0.336"
$ws.Range("E44").Value = "  +3.70%  "

# Row 45
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "3.32"
$ws.Range("E45").Value = "  -0.62%  "

# Row 46
$ws.Range("D46").Value = "0.0407"
$ws.Range("E46").Value = "  +0.90%  "

# Row 47
$ws.Range("D47").Value = "2.66"
$ws.Range("E47").Value = "  +7.36%  "

# Row 48
$ws.Range("E48").Value = "  +1.97%  "

# Row 49
$ws.Range("E49").Value = "  -0.75%  "

# Row 50
$ws.Range("E50").Value = "  -2.18%  "

# Row 51
$ws.Range("D51").Value = "129.81"
$ws.Range("E51").Value = "  +1.04%  "
